$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) from an existing header cell (G1) into the new header cell H1
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add the "Save" column data values
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
